# "Add cantrals by cantons"
#
# The sheet held two header rows (row1 = wide labels spanning merged-like
# columns, row2 = unit labels) followed by the data table, plus a stray
# "totals only" row part-way through the data (old row 22: just F/G/H
# totals, no idx/name/dates). The new layout collapses everything onto a
# single header row with explicit column titles (idx, idx2, Name, Date
# Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer,
# (GWh) Year) and removes the stray subtotal row, so every data row moves
# up accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stray subtotal-only row first (higher row index so the row2
# delete below doesn't shift its position), then collapse the old
# two-row header down to one row by deleting the old units row (row 2).
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(2).Delete()

# Rewrite the (now single) header row with the new column titles.
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 use the smaller (9pt) header font, same as the rest of the sheet's
# labels, but keep an explicit "General" number format of their own.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").NumberFormat = "General"

# Match the author's final selection (row 21 highlighted) left in the file.
$null = $ws.Range("A21:K21").Select()
